# Simulated Wild Card round and logged it
#
# Updates cumulative Rushing / Receiving stats for the Raiders roster with
# the results of the Wild Card game, and adds a new Receiving entry for
# D.Waller (who had not recorded a reception before this game).

$wb = $excel.ActiveWorkbook

$wsRushing   = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---------------------------------------------
# Row 2: D.Carr
$wsRushing.Cells.Item(2,5).Value = 13

# Row 3: M.Mariota
$wsRushing.Cells.Item(3,3).Value = 4
$wsRushing.Cells.Item(3,4).Value = 5
$wsRushing.Cells.Item(3,5).Value = 4
$wsRushing.Cells.Item(3,6).Value = 5

# Row 5: J.Jacobs
$wsRushing.Cells.Item(5,3).Value = 125
$wsRushing.Cells.Item(5,4).Value = 83
$wsRushing.Cells.Item(5,5).Value = 12
$wsRushing.Cells.Item(5,6).Value = 36

# Row 8: J.Richard
$wsRushing.Cells.Item(8,5).Value = 7
$wsRushing.Cells.Item(8,6).Value = 2

# Row 11: H.Renfrow
$wsRushing.Cells.Item(11,4).Value = 1

# --- Receiving sheet updates --------------------------------------------
# Row 2: J.Jacobs
$wsReceiving.Cells.Item(2,3).Value = 65
$wsReceiving.Cells.Item(2,4).Value = 55
$wsReceiving.Cells.Item(2,7).Value = 7
$wsReceiving.Cells.Item(2,8).Value = 4

# Row 5: J.Richard
$wsReceiving.Cells.Item(5,3).Value = 16
$wsReceiving.Cells.Item(5,4).Value = 13
$wsReceiving.Cells.Item(5,7).Value = 3

# Row 7: B.Edwards
$wsReceiving.Cells.Item(7,3).Value = 36
$wsReceiving.Cells.Item(7,4).Value = 24
$wsReceiving.Cells.Item(7,5).Value = 26
$wsReceiving.Cells.Item(7,6).Value = 11
$wsReceiving.Cells.Item(7,7).Value = 7

# Row 8: H.Renfrow
$wsReceiving.Cells.Item(8,3).Value = 120
$wsReceiving.Cells.Item(8,4).Value = 100
$wsReceiving.Cells.Item(8,5).Value = 18
$wsReceiving.Cells.Item(8,7).Value = 24
$wsReceiving.Cells.Item(8,8).Value = 17

# Row 9: Z.Jones
$wsReceiving.Cells.Item(9,3).Value = 51
$wsReceiving.Cells.Item(9,4).Value = 39
$wsReceiving.Cells.Item(9,5).Value = 22
$wsReceiving.Cells.Item(9,7).Value = 10
$wsReceiving.Cells.Item(9,8).Value = 4

# Row 11: D.Jackson
$wsReceiving.Cells.Item(11,3).Value = 13
$wsReceiving.Cells.Item(11,5).Value = 8
$wsReceiving.Cells.Item(11,6).Value = 4
$wsReceiving.Cells.Item(11,7).Value = 1

# Row 12: F.Moreau
$wsReceiving.Cells.Item(12,3).Value = 98
$wsReceiving.Cells.Item(12,4).Value = 76
$wsReceiving.Cells.Item(12,5).Value = 32
$wsReceiving.Cells.Item(12,6).Value = 26
$wsReceiving.Cells.Item(12,7).Value = 18
$wsReceiving.Cells.Item(12,8).Value = 11

# Row 13: D.Carrier
$wsReceiving.Cells.Item(13,3).Value = 23
$wsReceiving.Cells.Item(13,4).Value = 18
$wsReceiving.Cells.Item(13,5).Value = 3
$wsReceiving.Cells.Item(13,6).Value = 2
$wsReceiving.Cells.Item(13,7).Value = 4
$wsReceiving.Cells.Item(13,8).Value = 1

# Row 14: D.Helm
$wsReceiving.Cells.Item(14,3).Value = 2
$wsReceiving.Cells.Item(14,4).Value = 2

# Row 15 (new): D.Waller - first logged target/catch of the season
# Column A uses the same bold/bordered/centered style as the rest of the
# roster index column, so copy that formatting down from the row above.
$wsReceiving.Cells.Item(14,1).Copy()
$wsReceiving.Cells.Item(15,1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$wsReceiving.Cells.Item(15,1).Value = 13
$wsReceiving.Cells.Item(15,2).Value = "D.Waller"
$wsReceiving.Cells.Item(15,3).Value = 1
$wsReceiving.Cells.Item(15,4).Value = 1
$wsReceiving.Cells.Item(15,5).Value = 0
$wsReceiving.Cells.Item(15,6).Value = 0
$wsReceiving.Cells.Item(15,7).Value = 0
$wsReceiving.Cells.Item(15,8).Value = 0
